$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.100031495094299
$ws.Range("B1").Value = 1.940863251686096
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.014314413070679
$ws.Range("E1").Value = 1.124070763587952
